# The backing data source for this sheet was re-synced: every data row
# (rows 2-18, i.e. one header row followed by 17 records) now shows the
# content that used to live in a different row. The remap is a single
# 17-cycle permutation (no rows added or removed, row 1 header and rows
# 19-20 are untouched).
#
# destination row -> source row (content to copy into the destination)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 4
    3  = 5
    4  = 6
    5  = 7
    6  = 8
    7  = 9
    8  = 10
    9  = 11
    10 = 12
    11 = 2
    12 = 13
    13 = 14
    14 = 15
    15 = 16
    16 = 17
    17 = 18
    18 = 3
}

$firstCol = 1   # A
$lastCol  = 51  # AY

# 1) Snapshot every involved row's cell values (column by column) BEFORE any
#    writes happen. The mapping is a single permutation cycle, so every
#    source row is also a destination row that will later be overwritten -
#    we must capture everything up front.
$snapshots = @{}
foreach ($row in $mapping.Keys) {
    $rowVals = @{}
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $rowVals[$col] = $ws.Cells.Item($row, $col).Value2
    }
    $snapshots[$row] = $rowVals
}

# 2) Write each destination row from its snapshotted source row, cell by
#    cell. String values are written with the cell temporarily forced to
#    Text ("@") number format (restored to General right after), so that
#    text which happens to look like a date/time/number (e.g.
#    "2017-10-16") is written back verbatim as text instead of being
#    auto-converted by the usual Excel date/number literal parsing that
#    `.Value2` assignment triggers. Non-string values (numbers, booleans,
#    blanks) are assigned directly - they are already the right type and
#    don't need the Text-format workaround (all source cells use General
#    formatting, so this round-trips formatting exactly either way).
foreach ($row in $mapping.Keys) {
    $srcRow = $mapping[$row]
    $rowVals = $snapshots[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $val = $rowVals[$col]
        if ($val -is [string]) {
            $cell.NumberFormat = "@"
            $cell.Value2 = $val
            $cell.NumberFormat = "General"
        } else {
            $cell.Value2 = $val
        }
    }
}

Write-Output "Row remap applied."
